$wb = $excel.ActiveWorkbook

# --- 1) Status text: "Ready for handoff" -> "In Translation" ---
# This literal string is shared by the Overview sheet (columns E/F, row 2)
# and by each per-language sheet's "Status" column (column C, row 2).
# Updating every cell that held the old text keeps them all pointing at one
# (new) shared-string entry, mirroring the original single shared-string edit.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# --- 2) Narrow the "zh-cn"/"de-de" status columns ---
# Target stored column width is 13.4101845877511 characters. The ColumnWidth
# COM property is quantized to the nearest 1/6 of a character (Excel's
# pixel-grid rounding for the default Calibri 11 font), so 12.5 is the input
# that lands on the closest achievable grid point (13.333333333333334),
# nearest to the authored target width.
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F

$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C
